$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Insert a new paragraph "2015 - Present" right after
#    "Selection of Events:" (and before the blank bold paragraph that
#    used to precede "Selection of Keywords for Search Terms:").
#    The "_GoBack" bookmark that used to sit at the end of the
#    "To be continued..." paragraph moves here instead.
# ------------------------------------------------------------------

$found = $d.Content.Find.Execute("Selection of Events:", $false, $false,
    $false, $false, $false, $true, 1, $false, "", 0)

$selEvents = $d.Content
$anchorRange = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq "Selection of Events:") {
        $anchorRange = $para.Range
        break
    }
}

$anchorRange.InsertParagraphAfter()

# The freshly-created paragraph is the one right after "Selection of Events:"
$newParaIndex = $i + 1
$newPara = $d.Paragraphs.Item($newParaIndex)
$newRange = $newPara.Range
$newRange.Italic = $false
$newRange.Bold = $false
$newRange.InsertBefore("2015 - Present")

# ------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark onto the new paragraph (right after
#    its text, before the paragraph mark) and remove it from the
#    "To be continued..." paragraph.
# ------------------------------------------------------------------

try {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
} catch {
}

$newTextRange = $d.Paragraphs.Item($newParaIndex).Range.Duplicate
$newTextRange.Collapse(0)
$newTextRange.MoveEnd(1, -1)
$d.Bookmarks.Add("_GoBack", $newTextRange)

# ------------------------------------------------------------------
# 3. Tidy up a few citation runs whose stray spell-check markers
#    (proofErr) disappear once the surrounding text is touched again.
# ------------------------------------------------------------------

$d.Content.Find.Execute("B. Hague and B. Loader (eds), Digital Democracy",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "B. Hague and B. Loader (eds), Digital Democracy", 2)

$d.Content.Find.Execute("Heil, B. and Piskorski",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Heil, B. and Piskorski", 2)

$d.Content.Find.Execute("Enli, G. (2017)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Enli, G. (2017)", 2)
